$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 69

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 4).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2024-01-17"
$ws.Cells.Item($row, 2).Value = "22:43:12"
$ws.Cells.Item($row, 3).Value = "Wednesday"
$ws.Cells.Item($row, 4).Value = "02"

$ws.Cells.Item($row, 5).Value = 139005
$ws.Cells.Item($row, 6).Value = 139596
$ws.Cells.Item($row, 7).Value = 171266
$ws.Cells.Item($row, 8).Value = 148622
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 119991
$ws.Cells.Item($row, 11).Value = 222888
$ws.Cells.Item($row, 12).Value = 255058
$ws.Cells.Item($row, 13).Value = 185065
$ws.Cells.Item($row, 14).Value = 110359
$ws.Cells.Item($row, 15).Value = 41353
$ws.Cells.Item($row, 16).Value = 30937
$ws.Cells.Item($row, 17).Value = 73477
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42762
$ws.Cells.Item($row, 20).Value = -1
